# Generate Report for Handback
# - Update "Ready for handoff" -> "Handed back: in sync with en-US" on all three sheets
# - Fill in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#   columns on the zh-cn and de-de sheets now that the handback is in sync, including a
#   new hyperlink on the "Latest Target File" cell pointing at the source .md file
# - Widen a few columns that now hold longer values

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status columns E (zh-cn) and F (de-de) for both data rows
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the two status columns - they now hold the longer "Handed back..." text
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): status + handback columns
# ---------------------------------------------------------------------------
$mdUrl441f = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a6eeb829efe90fdeaae4a64f62e0b06a764c3ec/e2e/441f0192-99da-476c-b3df-fcf295381792.md"
$mdUrlCe8a = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a6eeb829efe90fdeaae4a64f62e0b06a764c3ec/e2e/ce8a05f4-787c-4f88-833a-79e42b51fa6c.md"

$langSheets = @{
    "zh-cn" = @{
        HandbackDateTimeRow2 = "2016-08-18 06:28:46"
        HandbackDateTimeRow3 = "2016-08-18 06:28:46"
        TargetFileRow2       = "441f0192-99da-476c-b3df-fcf295381792.43c4a393846b4c2b2669967eac3da9f49ea5e437.zh-cn.xlf"
        TargetFileRow3       = "ce8a05f4-787c-4f88-833a-79e42b51fa6c.8d99c69964948b7dc41341facfb3f790e5633eaa.zh-cn.xlf"
    }
    "de-de" = @{
        HandbackDateTimeRow2 = "2016-08-18 06:28:54"
        HandbackDateTimeRow3 = "2016-08-18 06:28:54"
        TargetFileRow2       = "441f0192-99da-476c-b3df-fcf295381792.43c4a393846b4c2b2669967eac3da9f49ea5e437.de-de.xlf"
        TargetFileRow3       = "ce8a05f4-787c-4f88-833a-79e42b51fa6c.8d99c69964948b7dc41341facfb3f790e5633eaa.de-de.xlf"
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $langSheets[$sheetName]

    # Status column (C) for both rows
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Row 2 (441f0192...): the file is now handed back and in sync, so the
    # "Latest Target File" (I) mirrors the source file name (and becomes a
    # hyperlink to the source .md), "Latest Handback File" (J) is the
    # generated xlf, and "Latest Handback DateTime" (K) is now a real timestamp
    $ws.Range("I2").Value = "441f0192-99da-476c-b3df-fcf295381792.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl441f, [Type]::Missing, [Type]::Missing, "441f0192-99da-476c-b3df-fcf295381792.md") | Out-Null
    $ws.Range("J2").Value = $info.TargetFileRow2
    $ws.Range("K2").Value = $info.HandbackDateTimeRow2

    # Row 3 (ce8a05f4...): same treatment
    $ws.Range("I3").Value = "ce8a05f4-787c-4f88-833a-79e42b51fa6c.md"
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrlCe8a, [Type]::Missing, [Type]::Missing, "ce8a05f4-787c-4f88-833a-79e42b51fa6c.md") | Out-Null
    $ws.Range("J3").Value = $info.TargetFileRow3
    $ws.Range("K3").Value = $info.HandbackDateTimeRow3

    # Widen columns that now contain longer values: Status (C), Latest Target
    # File (I) and Latest Handback File (J)
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}
